$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 41548
$ws.Range("A2").NumberFormat = "DD/MM/YY"

$ws.Range("B2").Value = 0.0833333333333333
$ws.Range("B2").NumberFormat = "HH:MM:SS"

$ws.Range("A3").Select() | Out-Null
